$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  "B" = @(8.936014356752448; 8.5634517970842; 8.327146427361857; 8.229090050531916; 8.212705956464303; 8.325830935086254; 8.809199737762382; 9.691932423644447; 10.29481393373; 10.55818291148871; 10.65628235904914; 10.63522846765248; 10.56628668345236; 10.52384340062427; 10.2773771196775; 10.12333886176422; 10.03371892567439; 10.00320198771016; 10.13984270918383; 10.58658135401009; 10.86901006672598; 10.71916524691582; 10.13238462021947; 9.460773792683158)
  "C" = @(5.822997384038935; 5.567647760220535; 5.403620833193349; 5.335010840440932; 5.323513157659692; 5.402702615761214; 5.73648514769413; 6.331634914449539; 6.730563384020052; 6.903392463923494; 6.96757359969664; 6.953807638903979; 6.908698160176597; 6.880901905911705; 6.719092630259647; 6.617596396422444; 6.558405743543526; 6.538225969186577; 6.628485129812839; 6.921982410150138; 7.106416010773396; 7.008662297335817; 6.623564945442931; 6.177232378862121)
  "D" = @(4.69325500866014; 4.638264397795136; 4.603726679817777; 4.589465623672242; 4.587086543557889; 4.603535095592431; 4.674459349130019; 4.807074921353851; 4.900138274903615; 4.941440531158698; 4.956925200635744; 4.953597330320538; 4.942717622320703; 4.936033017639938; 4.897417681597536; 4.873458784065241; 4.859581461735556; 4.854866434248412; 4.876019320084162; 4.94591753436077; 4.990690015025917; 4.966879640749268; 4.874862022114566; 4.771936155038496)
  "E" = @(16.47667668692905; 15.5442172489596; 14.94693422056099; 14.6975833398294; 14.65582819829871; 14.9435950957804; 16.16044781658172; 18.41992924820923; 20.05566072055843; 20.75765700656305; 21.01746595111139; 20.96177874503634; 20.7791521427036; 20.66650483313114; 20.00893915501725; 19.59478306624211; 19.35260624838605; 19.26992730444488; 19.63928083139394; 20.83295717116989; 21.57800907555512; 21.18356025048128; 19.61917608125287; 17.77999807161037)
  "F" = @(23.15625613623854; 23.18289448938697; 23.20656156170801; 23.21803850255876; 23.22005471693912; 23.20670893313944; 23.16392035400795; 23.13824139689659; 23.15510123277694; 23.17055239400391; 23.17752211719499; 23.17597132618081; 23.1711033830428; 23.16826728518236; 23.15424805405296; 23.14764096899667; 23.14457328220695; 23.14366041652268; 23.14826848775427; 23.17250286328959; 23.19486103952384; 23.18233194334196; 23.14798251000796; 23.13892608037337)
  "H" = @(7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261; 7.344005520526261)
  "I" = @(19.6611083466157; 19.76091451388519; 19.82664799666077; 19.85455235888398; 19.85925328745865; 19.82701980295097; 19.69459622013085; 19.47033349738275; 19.32729116686438; 19.26696396349264; 19.24480437313424; 19.24954632428772; 19.26512713510412; 19.27476012784855; 19.33132938920166; 19.36724976252177; 19.38835651086774; 19.39557945816313; 19.36337976623182; 19.26053205874485; 19.19730930848588; 19.23068603544623; 19.36512797238947; 19.52719734360459)
  "K" = @(8.689240716366754; 8.418152913604663; 8.248216453697529; 8.178193129835098; 8.166522331598378; 8.247275080429475; 8.596553438773929; 9.307098448307642; 9.874003388772795; 10.11979037569407; 10.21109032953148; 10.19150664301201; 10.12733741125112; 10.08779993137678; 9.857694075128069; 9.713406042121685; 9.629278020698399; 9.600599513772842; 9.728883728215918; 10.14623386586786; 10.4086427676699; 10.26954735900883; 9.721889923330989; 9.087862885637787)
  "O" = @(20.72629569191747; 20.80098525690403; 20.85225484366284; 20.87450258663557; 20.87827848372864; 20.8525494055709; 20.75092263419588; 20.59479282944617; 20.50670758719942; 20.47247939673865; 20.460363036455; 20.4629348475804; 20.47146561132649; 20.47680115365668; 20.50906247384037; 20.53035371951638; 20.5431494936496; 20.54757616313725; 20.52803031014968; 20.46893694436726; 20.43524413119111; 20.45277410366875; 20.52907899455659; 20.63237658057401)
}

$rows = 2..25
foreach ($col in $data.Keys) {
  $vals = $data[$col]
  for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    $ws.Range("$col$row").Value = $vals[$i]
  }
}